# daily auto push: 2026-01-17 06:45 UTC
# Insert a new daily log entry (2026/01/17, 13:xx) at row 641, pushing the
# existing rows (formerly 641-682) down by one to 642-683.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 641 (shifts 641..682 -> 642..683)
$ws.Rows.Item(641).Insert()

# The Date column stores plain text like "2026/01/17" (not a real date value)
# in every other row of this sheet, so force the cell to Text format before
# assigning the value - otherwise Excel will auto-convert the string into a
# date serial number.
$ws.Range("A641").NumberFormat = "@"
$ws.Range("A641").Value = "2026/01/17"
# Drop the temporary formatting so the new cell has no style override, just
# like the rest of the data cells in this sheet.
$ws.Range("A641:D641").ClearFormats()

$ws.Range("B641").Value = "土"
$ws.Range("C641").Value = 13
$ws.Range("D641").Value = 26
